# Apply profile/editor metadata updates across rows 2-18 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$editorImg    = "https://github.com/shoot649854/IMG_DB/blob/main/profile.webp"
$editorBio    = "Computer Science - University of California Santa Cruz | Intern - LiNK"
$editorSocial = "[{'icon': 'fab fa-linkedin', 'url': 'https://www.linkedin.com/in/shoto-morisaki-93b0a71bb/'}, {'icon': 'fab fa-github', 'url': 'https://github.com/shoot649854/'}, {'icon': 'fab fa-portfolio', 'url': 'https://portfolio-shoto.vercel.app/'}]"

for ($row = 2; $row -le 18; $row++) {
    # Clear relatedDoc1..relatedDoc5 (columns L..P), removing their "None" text entirely.
    # NOTE: use ${row} (braces) inside "L${row}:P${row}" -- "L$row:P$row" would be
    # mis-parsed by PowerShell as the variable "$row:P" (drive/scope syntax).
    $ws.Range("L${row}:P${row}").Value = ""

    # Update editor_img (Q) and editor_bio (R).
    $ws.Cells.Item($row, 17).Value = $editorImg
    $ws.Cells.Item($row, 18).Value = $editorBio

    # Update editor_social (T) with the new icon/url list.
    $ws.Cells.Item($row, 20).Value = $editorSocial
}

# Row 17 Status moves from Published to Draft.
$ws.Range("A17").Value = "Draft"
